$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
# Current Capital, Total P&L $/%, Total/Winning Trades, Win Rate % updated
# after trade #8 (MarketMaking) closed.
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B3").Value = 1199.9
$ws1.Range("B4").Value = -0.1
$ws1.Range("B5").Value = -0.25
$ws1.Range("B6").Value = 8
$ws1.Range("B7").Value = 3
$ws1.Range("B9").Value = 37.5

# ---- Sheet: Strategy Status ----
# MarketMaking strategy row (row 4) reflects the new trade too.
$ws2 = $wb.Worksheets.Item("Strategy Status")
$ws2.Range("C4").Value = 99.90000000000001
$ws2.Range("D4").Value = 8
$ws2.Range("E4").Value = -0.1
$ws2.Range("F4").Value = -0.1
$ws2.Range("G4").Value = 37.5

# ---- Sheet: All Trades & MarketMaking ----
# Append trade #8 (row 9) to both trade logs. Column B holds an ISO-style
# date string ("2026-02-17"); force text storage first (NumberFormat "@")
# so Excel doesn't auto-coerce it into a date serial number, then strip the
# temporary number-format back off so the cell keeps the workbook's default
# style, matching every other row in the sheet.
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A9").Value = 8

    $ws.Range("B9").NumberFormat = "@"
    $ws.Range("B9").Value = "2026-02-17"
    $ws.Range("B9").ClearFormats()

    $ws.Range("C9").Value = "15:14:05"
    $ws.Range("D9").Value = "MarketMaking"
    $ws.Range("E9").Value = "UP"
    $ws.Range("F9").Value = 0.9399999999999999
    $ws.Range("G9").Value = 0.95
    $ws.Range("H9").Value = "CLOSED"
    $ws.Range("I9").Value = 1.0638
    $ws.Range("J9").Value = 0.01
    $ws.Range("K9").Value = 99.90000000000001
    $ws.Range("L9").Value = 0
    $ws.Range("M9").Value = 0
    $ws.Range("N9").Value = 0.6
    $ws.Range("O9").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P9").Value = "early_exit"
    $ws.Range("Q9").Value = 0.13
}
